$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2
    if ($raw -eq $null) { continue }
    $val = [string]$raw
    if ($val -eq "") { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -le 1) { continue }

    $rev = $parts[-1..-($parts.Count)]
    $newval = [string]::Join(", ", $rev)

    $cell.Value = $newval
}
